$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values.
# Cells whose new text looks like a plain number must be pre-formatted as
# Text ("@") before assignment, otherwise Excel auto-converts the string
# into a numeric value (stripping formatting such as trailing zeros).
$ws.Range("D2").Value = "29.375.04"
$ws.Range("D3").Value = "1.842.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.31"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6271"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07396"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2893"
$ws.Range("D12").Value = "1.839.62"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.968"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6746"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001024"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.91"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.266"
$ws.Range("D18").Value = "29.365.69"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.37"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.311"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.71"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.486"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07285"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.475"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.477"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.040"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.818"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6999"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.570"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.894"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.796"
$ws.Range("D40").Value = "1.233.74"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9469"
$ws.Range("D43").Value = "1.990.83"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.97"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.26"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.968"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.843"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1128"

# Update Volume(1h) (column E) percentage-change values
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("E28").Value = "  +11.69%  "
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -2.24%  "

